$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8, col A was stored as text "76442781"; normalize it to a real number,
# matching the rest of the "phone" column.
$ws.Cells.Item(8, 1).Value = 76442781

# Append new payment row 9: 76442780 (Cash) 2025-08-15T09:54:44
# Phone numbers in this sheet are text-looking numerics that must stay text
# (see A8 pre-edit); build it as text via TEXT()+paste-values so Excel
# doesn't silently coerce the "quote-prefixed" numeric string into a number
# (and doesn't leave a stray NumberFormat style behind either).
$ws.Cells.Item(9, 1).Formula = '=TEXT(76442780,"0")'
$ws.Cells.Item(9, 1).Copy()
$ws.Cells.Item(9, 1).PasteSpecial(-4163)

$ws.Cells.Item(9, 2).Value = 17
$ws.Cells.Item(9, 3).Value = "Cash"
$ws.Cells.Item(9, 4).Value = "2025-08-15T09:54:44"
